# Fruta / hortaliza, semanal
# Insert a new weekly record at row 350, pushing existing records (old rows 350-439)
# down to rows 351-440.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 350; all rows 350:439 shift down to 351:440.
$ws.Rows("350:350").Insert()

# Populate the newly inserted row 350 with the new record's data.
$ws.Range("A350").Value = 5
$ws.Range("B350").Value = "Macroferia Regional de Talca"
$ws.Range("C350").Value = "Maule"
$ws.Range("D350").Value = 44855
$ws.Range("E350").Value = 7
$ws.Range("F350").Value = 100112032
$ws.Range("G350").Value = "Zapallo italiano"
$ws.Range("H350").Value = "Sin especificar"
$ws.Range("I350").Value = "Primera"
$ws.Range("J350").Value = 300
$ws.Range("K350").Value = 14000
$ws.Range("L350").Value = 14000
$ws.Range("M350").Value = 14000
$ws.Range("N350").Value = "`$/caja 50 unidades"
$ws.Range("O350").Value = "Región de O'Higgins"
$ws.Range("P350").Value = 280
$ws.Range("Q350").Value = 50
$ws.Range("R350").Value = "Hortaliza"
